$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "334.86"
$ws.Range("E2").Value = "1.61%"

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "43.89"
$ws.Range("E3").Value = "6.42%"

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.755"
$ws.Range("E4").Value = "1.98%"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08363"
$ws.Range("E5").Value = "1.89%"

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "8.848"
$ws.Range("E6").Value = "1.04%"

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "4.514"
$ws.Range("E7").Value = "0.24%"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.60%"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.89%"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9474"
$ws.Range("E10").Value = "2.62%"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1248"
$ws.Range("E11").Value = "-2.33%"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1971"
$ws.Range("E12").Value = "0.77%"

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1039"
$ws.Range("E13").Value = "11.32%"

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04558"
$ws.Range("E14").Value = "18.72%"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.73%"

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001304"
$ws.Range("E16").Value = "-0.68%"

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005909"
$ws.Range("E17").Value = "-4.94%"

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "3.497"
$ws.Range("E18").Value = "1.42%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.75%"

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "8.691"
$ws.Range("E20").Value = "4.20%"

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1363"
$ws.Range("E21").Value = "-0.20%"

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2638"
$ws.Range("E22").Value = "-0.84%"

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04412"
$ws.Range("E23").Value = "0.21%"

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001256"
$ws.Range("E24").Value = "-0.08%"

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004350"
$ws.Range("E25").Value = "0.79%"

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001261"
$ws.Range("E26").Value = "5.03%"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.94%"

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06056"
$ws.Range("E40").Value = "10.19%"

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007912"
$ws.Range("E41").Value = "0.39%"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1426"
$ws.Range("E42").Value = "0.36%"

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008977"
$ws.Range("E43").Value = "0.43%"

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002144"
$ws.Range("E44").Value = "-1.21%"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-11.15%"

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007269"
$ws.Range("E46").Value = "7.35%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.09%"

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003186"
$ws.Range("E48").Value = "-0.08%"

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002272"
$ws.Range("E49").Value = "-0.36%"

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "0.09%"

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "0.09%"
